$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("H2").Value = 0.003000259399414062

# Row 3
$ws.Range("C3").Value = 514.34
$ws.Range("D3").Value = 512.67
$ws.Range("E3").Value = 1.67
$ws.Range("F3").Value = 1027.01
$ws.Range("G3").Value = 513.5
$ws.Range("H3").Value = 1.598640203475952

# Row 4
$ws.Range("H4").Value = 0.003998994827270508

# Row 5
$ws.Range("C5").Value = 627.89
$ws.Range("D5").Value = 491.7
$ws.Range("E5").Value = 136.19
$ws.Range("F5").Value = 1119.59
$ws.Range("G5").Value = 559.8
$ws.Range("H5").Value = 1.647355318069458

# Row 6
$ws.Range("H6").Value = 0.007999897003173828

# Row 7
$ws.Range("C7").Value = 553.34
$ws.Range("D7").Value = 551.22
$ws.Range("E7").Value = 2.11
$ws.Range("F7").Value = 1656.81
$ws.Range("G7").Value = 552.27
$ws.Range("H7").Value = 3.590090274810791

# Row 8
$ws.Range("H8").Value = 0.008063077926635742

# Row 9
$ws.Range("C9").Value = 640.11
$ws.Range("D9").Value = 513.79
$ws.Range("E9").Value = 126.32
$ws.Range("F9").Value = 1792.56
$ws.Range("G9").Value = 597.52
$ws.Range("H9").Value = 3.861264705657959

# Row 10
$ws.Range("H10").Value = 0.01200008392333984

# Row 11
$ws.Range("C11").Value = 543.09
$ws.Range("D11").Value = 515.6900000000001
$ws.Range("E11").Value = 27.4
$ws.Range("F11").Value = 2139.9
$ws.Range("G11").Value = 534.97
$ws.Range("H11").Value = 6.09990930557251

# Row 12
$ws.Range("H12").Value = 0.01099920272827148

# Row 13
$ws.Range("C13").Value = 701.87
$ws.Range("D13").Value = 465.73
$ws.Range("E13").Value = 236.14
$ws.Range("F13").Value = 2332.26
$ws.Range("G13").Value = 583.0599999999999
$ws.Range("H13").Value = 6.773151874542236
